$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '69.515.37'
$ws.Cells.Item(2, 5).Value = '  +0.20%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.674.30'
$ws.Cells.Item(3, 5).Value = '  -0.38%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '645.91'
$ws.Cells.Item(5, 5).Value = '  -4.87%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '159.43'
$ws.Cells.Item(6, 5).Value = '  -0.77%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  -0.10%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.499'
$ws.Cells.Item(8, 5).Value = '  +1.01%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.145'
$ws.Cells.Item(9, 5).Value = '  -0.84%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.16'
$ws.Cells.Item(10, 5).Value = '  -0.27%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.442'
$ws.Cells.Item(11, 5).Value = '  +0.62%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000231'
$ws.Cells.Item(12, 5).Value = '  -0.87%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.288.01'
$ws.Cells.Item(13, 5).Value = '  -0.52%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '32.64'
$ws.Cells.Item(14, 5).Value = '  +0.57%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.684.84'
$ws.Cells.Item(15, 5).Value = '  -0.33%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '69.482.83'
$ws.Cells.Item(16, 5).Value = '  +0.18%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.117'
$ws.Cells.Item(17, 5).Value = '  +0.01%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '15.97'
$ws.Cells.Item(18, 5).Value = '  -0.34%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.48'
$ws.Cells.Item(19, 5).Value = '  +0.28%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '471.46'
$ws.Cells.Item(20, 5).Value = '  -0.16%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.06'
$ws.Cells.Item(21, 5).Value = '  +2.61%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.649'
$ws.Cells.Item(22, 5).Value = '  -0.17%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '79.42'
$ws.Cells.Item(23, 5).Value = '  -1.05%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.817.62'
$ws.Cells.Item(24, 5).Value = '  -0.46%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.07%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.0000124'
$ws.Cells.Item(26, 5).Value = '  -0.73%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.95'
$ws.Cells.Item(27, 5).Value = '  +0.46%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.04'
$ws.Cells.Item(28, 5).Value = '  -0.83%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.63'
$ws.Cells.Item(29, 5).Value = '  -2.85%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.70'
$ws.Cells.Item(30, 5).Value = '  -2.24%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.00'
$ws.Cells.Item(31, 5).Value = '  -0.60%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.999'
$ws.Cells.Item(32, 5).Value = '  -0.21%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '26.80'
$ws.Cells.Item(33, 5).Value = '  -0.65%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '6.44'
$ws.Cells.Item(34, 5).Value = '  -2.17%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.77%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.669.31'
$ws.Cells.Item(36, 5).Value = '  -0.24%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '8.45'
$ws.Cells.Item(37, 5).Value = '  -0.20%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.87'
$ws.Cells.Item(39, 5).Value = '  -5.49%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '178.41'
$ws.Cells.Item(40, 5).Value = '  +5.82%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.03%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.20'
$ws.Cells.Item(42, 5).Value = '  -2.51%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.0899'
$ws.Cells.Item(43, 5).Value = '  -0.45%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.926'
$ws.Cells.Item(44, 5).Value = '  -1.53%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '47.21'
$ws.Cells.Item(45, 5).Value = '  +1.15%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '28.92'
$ws.Cells.Item(46, 5).Value = '  +2.89%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.69'
$ws.Cells.Item(47, 5).Value = '  -1.22%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.07'
$ws.Cells.Item(48, 5).Value = '  -1.28%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'FLOKI'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.000266'
$ws.Cells.Item(49, 5).Value = '  -4.58%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Cosmos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.81'
$ws.Cells.Item(50, 5).Value = '  -0.88%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.23'
$ws.Cells.Item(51, 5).Value = '  -4.45%  '
